$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G (MSRP_VERSION) mirrors existing column F (VERSION):
#  - data rows carry the same "normal text" formatting as F2/F3
#  - header row carries the same "Good" header formatting as F1
# Cells are written data-rows-first, then the header, so the shared-string
# table grows in the same order the workbook records them.
$ws.Range("F2").Copy()
$ws.Range("G2").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G3").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G2").Value = "MSRP_2000"
$ws.Range("G3").Value = "MSRP_2000_SS_TEST"

$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("G1").Value = "MSRP_VERSION"

# Auto-size the new column like Excel does for a freshly added column
$ws.Columns.Item(7).AutoFit()
$excel.CutCopyMode = $false

# Match the selection left behind after the edit
$ws.Range("H8").Select()
